$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 30, 33, 34, 37, 38, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49, 51, 57, 60, 61, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76, 78, 86, 87, 88, 89, 90, 93, 95, 96, 97, 99, 102, 104, 112, 113, 114, 115, 116, 119, 121, 122, 123, 125, 128, 130, 138, 139, 140, 141, 142, 145, 147, 148, 149, 151, 154, 156)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $current = $cell.Value2
    $parts = $current -split ", "
    if ($parts.Count -eq 2) {
        $newValue = "$($parts[1]), $($parts[0])"
        $cell.Value = $newValue
    }
}
